$d = $word.ActiveDocument

# The "Schedule" table has a 3rd ("Due") column that is blank for several
# rows. Those blank cells hold a single, otherwise-unformatted paragraph
# (just the end-of-cell/paragraph marks, no real text run). Give each of
# those empty paragraphs the "Compact" style, matching the already-styled
# sibling paragraphs elsewhere in the table.
$table = $d.Tables.Item(1)

foreach ($row in $table.Rows) {
    foreach ($cell in $row.Cells) {
        foreach ($para in $cell.Range.Paragraphs) {
            $bodyText = $para.Range.Text.TrimEnd([char]13, [char]7)
            if ($bodyText.Length -eq 0) {
                $para.Style = "Compact"
            }
        }
    }
}
